# chore: adapt column header formatting to respective input file names
#
# - Rename the "_old" / "_new" column header suffixes to the concrete
#   format-version names "_FV2310" / "_FV2404".
# - Turn the used range into a proper Excel Table ("Table1").
# - Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the header cells in row 1 (columns A:J = "_old" -> "_FV2310",
#    column K is the unchanged "diff" column, columns L:U = "_new" -> "_FV2404").
$headerNames = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headerNames.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headerNames[$i]
}

# 2. Convert the populated range A1:U65 into an Excel Table so the header
#    row gets AutoFilter buttons and the rows become structured data.
$usedRange = $ws.Range("A1:U65")
$table = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# 3. Freeze the header row (split below row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
